$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the new rows as Text so numeric-looking strings
# (amounts, phone numbers, comma-separated id lists, etc.) are
# preserved exactly as text instead of being coerced to numbers.
$ws.Range("A61:BV63").NumberFormat = "@"

# Row 61
$ws.Range("A61").Value = 'SALESIAN EDUCATION SOCIETY'
$ws.Range("B61").Value = '753702'
$ws.Range("C61").Value = '1234'
$ws.Range("D61").Value = '11000316561861'
$ws.Range("E61").Value = '1763044851'
$ws.Range("F61").Value = '11250.00'
$ws.Range("G61").Value = 'INR'
$ws.Range("H61").Value = '13-Nov-2025 20:24:29'
$ws.Range("I61").Value = 'Multi'
$ws.Range("J61").Value = 'sale'
$ws.Range("K61").Value = 'ICICI UPI QR'
$ws.Range("L61").Value = 'OK'
$ws.Range("M61").Value = 'NRNS'
$ws.Range("N61").Value = '108566739973'
$ws.Range("O61").Value = ''
$ws.Range("P61").Value = '100000036600'
$ws.Range("Q61").Value = 'IFSC0000000'
$ws.Range("R61").Value = ''
$ws.Range("S61").Value = 'MERCHANT'
$ws.Range("T61").Value = 'UPI'
$ws.Range("U61").Value = ''
$ws.Range("V61").Value = ''
$ws.Range("W61").Value = ''
$ws.Range("X61").Value = 'JAHNAVI KOLASANI'
$ws.Range("Y61").Value = 'kotakschoolvsp@gmail.com'
$ws.Range("Z61").Value = '9347559040'
$ws.Range("AA61").Value = ''
$ws.Range("AB61").Value = ''
$ws.Range("AC61").Value = ''
$ws.Range("AD61").Value = ''
$ws.Range("AE61").Value = ''
$ws.Range("AF61").Value = ''
$ws.Range("AG61").Value = ''
$ws.Range("AH61").Value = ''
$ws.Range("AI61").Value = '0.00'
$ws.Range("AJ61").Value = '0.00'
$ws.Range("AK61").Value = ''
$ws.Range("AL61").Value = ''
$ws.Range("AM61").Value = ''
$ws.Range("AN61").Value = ''
$ws.Range("AO61").Value = ''
$ws.Range("AP61").Value = ''
$ws.Range("AQ61").Value = ''
$ws.Range("AR61").Value = ''
$ws.Range("AS61").Value = ''
$ws.Range("AT61").Value = 'TRANSACTION IS SUCCESSFUL'
$ws.Range("AU61").Value = ''
$ws.Range("AV61").Value = ''
$ws.Range("AW61").Value = ''
$ws.Range("AX61").Value = ''
$ws.Range("AY61").Value = ''
$ws.Range("AZ61").Value = ''
$ws.Range("BA61").Value = 'UPI'
$ws.Range("BB61").Value = ''
$ws.Range("BC61").Value = ''
$ws.Range("BD61").Value = 'REGULAR'
$ws.Range("BE61").Value = '19295'
$ws.Range("BF61").Value = '263081,264584'
$ws.Range("BG61").Value = '2037,2044'
$ws.Range("BH61").Value = 'eleven thousand two hundred fifty'
$ws.Range("BI61").Value = '15599'
$ws.Range("BJ61").Value = ''
$ws.Range("BK61").Value = ''
$ws.Range("BL61").Value = ''
$ws.Range("BM61").Value = ''
$ws.Range("BN61").Value = ''
$ws.Range("BO61").Value = ''
$ws.Range("BP61").Value = ''
$ws.Range("BQ61").Value = ''
$ws.Range("BR61").Value = ''
$ws.Range("BS61").Value = ''
$ws.Range("BT61").Value = '0'
$ws.Range("BU61").Value = 'NA'
$ws.Range("BV61").Value = '-'

# Row 62
$ws.Range("A62").Value = 'SALESIAN EDUCATION SOCIETY'
$ws.Range("B62").Value = '753702'
$ws.Range("C62").Value = '1234'
$ws.Range("D62").Value = '11000316561861'
$ws.Range("E62").Value = '1763044851'
$ws.Range("F62").Value = '10750.00'
$ws.Range("G62").Value = 'INR'
$ws.Range("H62").Value = '13-Nov-2025 20:24:29'
$ws.Range("I62").Value = 'VlllX'
$ws.Range("J62").Value = 'sale'
$ws.Range("K62").Value = 'ICICI UPI QR'
$ws.Range("L62").Value = 'OK'
$ws.Range("M62").Value = 'NRNS'
$ws.Range("N62").Value = '108566739973'
$ws.Range("O62").Value = ''
$ws.Range("P62").Value = '0899053000000003'
$ws.Range("Q62").Value = 'SIBL0000899'
$ws.Range("R62").Value = ''
$ws.Range("S62").Value = 'MERCHANT'
$ws.Range("T62").Value = 'UPI'
$ws.Range("U62").Value = ''
$ws.Range("V62").Value = ''
$ws.Range("W62").Value = ''
$ws.Range("X62").Value = 'JAHNAVI KOLASANI'
$ws.Range("Y62").Value = 'kotakschoolvsp@gmail.com'
$ws.Range("Z62").Value = '9347559040'
$ws.Range("AA62").Value = ''
$ws.Range("AB62").Value = ''
$ws.Range("AC62").Value = ''
$ws.Range("AD62").Value = ''
$ws.Range("AE62").Value = ''
$ws.Range("AF62").Value = ''
$ws.Range("AG62").Value = ''
$ws.Range("AH62").Value = ''
$ws.Range("AI62").Value = '0.00'
$ws.Range("AJ62").Value = '0.00'
$ws.Range("AK62").Value = ''
$ws.Range("AL62").Value = ''
$ws.Range("AM62").Value = ''
$ws.Range("AN62").Value = ''
$ws.Range("AO62").Value = ''
$ws.Range("AP62").Value = ''
$ws.Range("AQ62").Value = ''
$ws.Range("AR62").Value = ''
$ws.Range("AS62").Value = ''
$ws.Range("AT62").Value = 'TRANSACTION IS SUCCESSFUL'
$ws.Range("AU62").Value = ''
$ws.Range("AV62").Value = ''
$ws.Range("AW62").Value = ''
$ws.Range("AX62").Value = ''
$ws.Range("AY62").Value = ''
$ws.Range("AZ62").Value = ''
$ws.Range("BA62").Value = 'UPI'
$ws.Range("BB62").Value = ''
$ws.Range("BC62").Value = ''
$ws.Range("BD62").Value = 'REGULAR'
$ws.Range("BE62").Value = '19295'
$ws.Range("BF62").Value = '263081,264584'
$ws.Range("BG62").Value = '2037,2044'
$ws.Range("BH62").Value = 'eleven thousand two hundred fifty'
$ws.Range("BI62").Value = '15599'
$ws.Range("BJ62").Value = ''
$ws.Range("BK62").Value = ''
$ws.Range("BL62").Value = ''
$ws.Range("BM62").Value = ''
$ws.Range("BN62").Value = ''
$ws.Range("BO62").Value = ''
$ws.Range("BP62").Value = ''
$ws.Range("BQ62").Value = ''
$ws.Range("BR62").Value = ''
$ws.Range("BS62").Value = ''
$ws.Range("BT62").Value = '0'
$ws.Range("BU62").Value = 'NA'
$ws.Range("BV62").Value = '-'

# Row 63
$ws.Range("A63").Value = 'SALESIAN EDUCATION SOCIETY'
$ws.Range("B63").Value = '753702'
$ws.Range("C63").Value = '1234'
$ws.Range("D63").Value = '11000316561861'
$ws.Range("E63").Value = '1763044851'
$ws.Range("F63").Value = '500.00'
$ws.Range("G63").Value = 'INR'
$ws.Range("H63").Value = '13-Nov-2025 20:24:29'
$ws.Range("I63").Value = 'lVl'
$ws.Range("J63").Value = 'sale'
$ws.Range("K63").Value = 'ICICI UPI QR'
$ws.Range("L63").Value = 'OK'
$ws.Range("M63").Value = 'NRNS'
$ws.Range("N63").Value = '108566739973'
$ws.Range("O63").Value = ''
$ws.Range("P63").Value = '0899053000000002'
$ws.Range("Q63").Value = 'SIBL0000899'
$ws.Range("R63").Value = ''
$ws.Range("S63").Value = 'MERCHANT'
$ws.Range("T63").Value = 'UPI'
$ws.Range("U63").Value = ''
$ws.Range("V63").Value = ''
$ws.Range("W63").Value = ''
$ws.Range("X63").Value = 'JAHNAVI KOLASANI'
$ws.Range("Y63").Value = 'kotakschoolvsp@gmail.com'
$ws.Range("Z63").Value = '9347559040'
$ws.Range("AA63").Value = ''
$ws.Range("AB63").Value = ''
$ws.Range("AC63").Value = ''
$ws.Range("AD63").Value = ''
$ws.Range("AE63").Value = ''
$ws.Range("AF63").Value = ''
$ws.Range("AG63").Value = ''
$ws.Range("AH63").Value = ''
$ws.Range("AI63").Value = '0.00'
$ws.Range("AJ63").Value = '0.00'
$ws.Range("AK63").Value = ''
$ws.Range("AL63").Value = ''
$ws.Range("AM63").Value = ''
$ws.Range("AN63").Value = ''
$ws.Range("AO63").Value = ''
$ws.Range("AP63").Value = ''
$ws.Range("AQ63").Value = ''
$ws.Range("AR63").Value = ''
$ws.Range("AS63").Value = ''
$ws.Range("AT63").Value = 'TRANSACTION IS SUCCESSFUL'
$ws.Range("AU63").Value = ''
$ws.Range("AV63").Value = ''
$ws.Range("AW63").Value = ''
$ws.Range("AX63").Value = ''
$ws.Range("AY63").Value = ''
$ws.Range("AZ63").Value = ''
$ws.Range("BA63").Value = 'UPI'
$ws.Range("BB63").Value = ''
$ws.Range("BC63").Value = ''
$ws.Range("BD63").Value = 'REGULAR'
$ws.Range("BE63").Value = '19295'
$ws.Range("BF63").Value = '263081,264584'
$ws.Range("BG63").Value = '2037,2044'
$ws.Range("BH63").Value = 'eleven thousand two hundred fifty'
$ws.Range("BI63").Value = '15599'
$ws.Range("BJ63").Value = ''
$ws.Range("BK63").Value = ''
$ws.Range("BL63").Value = ''
$ws.Range("BM63").Value = ''
$ws.Range("BN63").Value = ''
$ws.Range("BO63").Value = ''
$ws.Range("BP63").Value = ''
$ws.Range("BQ63").Value = ''
$ws.Range("BR63").Value = ''
$ws.Range("BS63").Value = ''
$ws.Range("BT63").Value = '0'
$ws.Range("BU63").Value = 'NA'
$ws.Range("BV63").Value = '-'

Write-Host "Rows 61-63 added"